$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555290451293.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555290701232.csv"
$ws1.Range("B4").Value = "go_stims-16512555290721257.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555290861254.csv"
$ws1.Name = "GNG_TO-16512555290871246"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_7-16512555297424881.csv"
$ws2.Range("B3").Value = "TB-16512555314916477.csv"
$ws2.Range("B4").Value = "ZB-match_6-16512555291231256.csv"
$ws2.Range("B5").Value = "OB-16512555299154882.csv"
$ws2.Range("B6").Value = "OB-1651255529854487.csv"
$ws2.Range("B7").Value = "TB-16512555313256469.csv"
$ws2.Range("B8").Value = "ZB-match_6-16512555294011624.csv"
$ws2.Range("B9").Value = "TB-16512555308456485.csv"
$ws2.Range("B10").Value = "OB-1651255530095489.csv"
$ws2.Name = "NB_TO-16512555315016506"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"
$ws3.Name = "RS_TO-16512555315086603"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555315316484.csv"
$ws4.Range("B3").Value = "ZM_stims-1651255531510648.csv"
$ws4.Range("B4").Value = "MM_stims-1651255531547647.csv"
$ws4.Range("B5").Value = "ZM_stims-1651255531532648.csv"
$ws4.Range("B6").Value = "MM_stims-16512555315636466.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555315486474.csv"
$ws4.Name = "TOL_TO-1651255531564647"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512555316276498.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555315716531.csv"
$ws5.Range("B4").Value = "vSAT_stims-16512555316106465.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555315946486.csv"
$ws5.Name = "vSAT_TO-16512555316426508"
